$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously empty PriceChange / UpDown for row 7 ---
$ws.Range("X7").Value = 0.42000000000000171
$ws.Range("Y7").Value = "Up"

# --- Append a brand new data row (row 8) ---
$ws.Range("A8").Value = 42649.879895833335
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 12
$ws.Range("E8").Value = 4014
$ws.Range("F8").Value = 683
$ws.Range("G8").Value = 51
$ws.Range("H8").Value = 48
$ws.Range("I8").Value = 70
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 7287
$ws.Range("L8").Value = 83
$ws.Range("M8").Value = 78
$ws.Range("N8").Value = 35
$ws.Range("O8").Value = 15
$ws.Range("P8").Value = "Noun"
$ws.Range("Q8").Value = 47.321424984051369
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.0834
$ws.Range("T8").Value = -0.0062
$ws.Range("U8").Value = 2.31
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0

# Match number formats of the columns above (date / percentage styles) by
# copying formats only, so the existing style indexes are reused instead of
# new (duplicate) number formats being created.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("S7").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("T7").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column C got a touch wider to fit "Neutral" ---
$ws.Columns.Item(3).ColumnWidth = 5.71
